# Refresh the cryptocurrency price/volume snapshot on the worksheet.
#
# The sheet stores every data cell (price in column D, 1h volume-change in
# column E, plus the coin name/link in B/C for a couple of reordered rows)
# as literal text - even values that look like plain numbers, e.g. "4.20" -
# so that things like "2.433.57" (thousand-separated, not a valid float) or
# "0.0₃0784" (subscript digit-grouping notation) round-trip untouched.
#
# Assigning a numeric-looking string straight to Range.Value makes Excel
# "smart" and silently reinterprets it as a real number (losing the text
# type and any formatting quirks). To keep those cells as text we prefix
# the assignment with an apostrophe (Excel's own "force text" quote
# prefix) and then reset the cell Style back to Normal so the temporary
# quote-prefix marker doesn't linger as a spurious style change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '61.956.81' },
    @{ Cell = 'E2'; Value = '  +3.03%  ' },
    @{ Cell = 'D3'; Value = '2.431.11' },
    @{ Cell = 'E3'; Value = '  +4.82%  ' },
    @{ Cell = 'E4'; Value = '  -0.06%  ' },
    @{ Cell = 'D5'; Value = '554.02' },
    @{ Cell = 'E5'; Value = '  +2.03%  ' },
    @{ Cell = 'D6'; Value = '138.51' },
    @{ Cell = 'E6'; Value = '  +6.01%  ' },
    @{ Cell = 'E7'; Value = '  -0.03%  ' },
    @{ Cell = 'D8'; Value = '0.582' },
    @{ Cell = 'E8'; Value = '  +1.00%  ' },
    @{ Cell = 'D9'; Value = '2.429.39' },
    @{ Cell = 'E9'; Value = '  +4.83%  ' },
    @{ Cell = 'E10'; Value = '  +3.20%  ' },
    @{ Cell = 'D11'; Value = '5.71' },
    @{ Cell = 'E11'; Value = '  +3.39%  ' },
    @{ Cell = 'E12'; Value = '  +0.28%  ' },
    @{ Cell = 'E13'; Value = '  +4.48%  ' },
    @{ Cell = 'D14'; Value = '26.02' },
    @{ Cell = 'E14'; Value = '  +10.96%  ' },
    @{ Cell = 'D15'; Value = '2.864.80' },
    @{ Cell = 'E15'; Value = '  +4.76%  ' },
    @{ Cell = 'D16'; Value = '61.866.82' },
    @{ Cell = 'E16'; Value = '  +2.90%  ' },
    @{ Cell = 'D17'; Value = '0.0000141' },
    @{ Cell = 'D18'; Value = '2.432.38' },
    @{ Cell = 'E18'; Value = '  +5.04%  ' },
    @{ Cell = 'E19'; Value = '  +5.88%  ' },
    @{ Cell = 'D20'; Value = '343.02' },
    @{ Cell = 'E20'; Value = '  +9.55%  ' },
    @{ Cell = 'D21'; Value = '4.18' },
    @{ Cell = 'E21'; Value = '  +1.91%  ' },
    @{ Cell = 'D22'; Value = '6.79' },
    @{ Cell = 'E22'; Value = '  +2.90%  ' },
    @{ Cell = 'E23'; Value = '  -0.05%  ' },
    @{ Cell = 'D24'; Value = '64.96' },
    @{ Cell = 'E24'; Value = '  +1.89%  ' },
    @{ Cell = 'E25'; Value = '  +0.62%  ' },
    @{ Cell = 'E26'; Value = '  +0.09%  ' },
    @{ Cell = 'D27'; Value = '1.51' },
    @{ Cell = 'E27'; Value = '  +12.47%  ' },
    @{ Cell = 'D28'; Value = '8.20' },
    @{ Cell = 'E28'; Value = '  +5.81%  ' },
    @{ Cell = 'E29'; Value = '  +11.75%  ' },
    @{ Cell = 'B30'; Value = 'PancakeSwap' },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' },
    @{ Cell = 'D30'; Value = '1.80' },
    @{ Cell = 'E30'; Value = '  +4.72%  ' },
    @{ Cell = 'B31'; Value = 'PEPE' },
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' },
    @{ Cell = 'D31'; Value = '0.0₃0783' },
    @{ Cell = 'E31'; Value = '  +7.61%  ' },
    @{ Cell = 'D32'; Value = '171.85' },
    @{ Cell = 'E32'; Value = '  +0.29%  ' },
    @{ Cell = 'D33'; Value = '6.30' },
    @{ Cell = 'E33'; Value = '  +7.30%  ' },
    @{ Cell = 'D34'; Value = '1.43' },
    @{ Cell = 'E34'; Value = '  +4.59%  ' },
    @{ Cell = 'D35'; Value = '0.394' },
    @{ Cell = 'E35'; Value = '  +3.94%  ' },
    @{ Cell = 'D36'; Value = '18.52' },
    @{ Cell = 'E36'; Value = '  +4.58%  ' },
    @{ Cell = 'E37'; Value = '  +11.04%  ' },
    @{ Cell = 'D38'; Value = '365.30' },
    @{ Cell = 'E38'; Value = '  +14.74%  ' },
    @{ Cell = 'E39'; Value = '  -0.02%  ' },
    @{ Cell = 'D40'; Value = '0.999' },
    @{ Cell = 'E40'; Value = '  -0.19%  ' },
    @{ Cell = 'D41'; Value = '1.69' },
    @{ Cell = 'E41'; Value = '  +11.04%  ' },
    @{ Cell = 'D42'; Value = '39.21' },
    @{ Cell = 'E42'; Value = '  +3.50%  ' },
    @{ Cell = 'D43'; Value = '145.98' },
    @{ Cell = 'E43'; Value = '  +6.84%  ' },
    @{ Cell = 'D44'; Value = '3.65' },
    @{ Cell = 'E44'; Value = '  +5.93%  ' },
    @{ Cell = 'D45'; Value = '20.61' },
    @{ Cell = 'E45'; Value = '  +9.35%  ' },
    @{ Cell = 'D46'; Value = '0.0953' },
    @{ Cell = 'E46'; Value = '  +1.43%  ' },
    @{ Cell = 'E47'; Value = '  +4.11%  ' },
    @{ Cell = 'E48'; Value = '  +4.87%  ' },
    @{ Cell = 'B49'; Value = 'VeChain' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D49'; Value = '0.0220' },
    @{ Cell = 'E49'; Value = '  +4.01%  ' },
    @{ Cell = 'B50'; Value = 'EnergySwap' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D50'; Value = '17.82' },
    @{ Cell = 'E50'; Value = '  +6.14%  ' },
    @{ Cell = 'E51'; Value = '  -6.45%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $value = $u.Value

    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        # Force text storage (quote prefix), then strip the resulting
        # "quote prefix" style back off so the cell keeps its original,
        # unstyled look.
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
